$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: CP002_login_exitoso gets credentials + welcome message
$wb.Worksheets.Item(1).Hyperlinks.Add($ws.Range("B3"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B3").Value = "jisola.tsoft@gmail.com"
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("C3").Value = 12061990
$ws.Range("D3").Value = "Te damos la bienvenida a Facebook, Juan"

# Row 4: CP003_cerrar_sesion gets credentials too
$wb.Worksheets.Item(1).Hyperlinks.Add($ws.Range("B4"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B4").Value = "jisola.tsoft@gmail.com"
$ws.Range("B4").Style = $ws.Range("B2").Style
$ws.Range("C4").Value = 12061990

# Update selection to D4
$ws.Range("D4").Select()
